# Work diary for the day - add new journal entries (rows 26-30) to the
# "Journal de travail" table on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: complete the existing date-only row with the rest of the entry
$ws.Range("B26").Value = "Conception"
$ws.Range("C26").Value = 2.25
$ws.Range("D26").Value = "Création des maquettes "
$ws.Range("E26").Value = "Home, login, register"

# --- Row 27
$ws.Range("A27").Value = 44687
$ws.Range("B27").Value = "Conception"
$ws.Range("C27").Value = 0.75
$ws.Range("D27").Value = "Révision du MLD"

# --- Row 28 (taller row - wrapped remark text)
$ws.Range("A28").Value = 44687
$ws.Range("B28").Value = "Review"
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = "Sprint review avec chef de projet "
$ws.Range("E28").Value = "Plusieurs choses à réviser. A voir dans Documentation/reviews"
$ws.Range("F28").Value = "601ce3acc8b6cbc18d49aa7c882488e928093d11"

# --- Row 29 (Commit ID filled in later, after row 30 - see below)
$ws.Range("A29").Value = 44687
$ws.Range("B29").Value = "Analyse"
$ws.Range("C29").Value = 1.25
$ws.Range("D29").Value = "Révision de la doc selon review"

# --- Row 30 (taller row - long wrapped remark)
$ws.Range("A30").Value = 44687
$ws.Range("B30").Value = "Analyse"
$ws.Range("C30").Value = 0.25
$ws.Range("D30").Value = "Préparation du second rendu"
$ws.Range("E30").Value = "Il reste la stratégie de test a finaliser dans la documentation. Le reste est plutôt bon et je pense être prêt à commencer la mise en place du site la semaine prochaine "

# The commit id for the review-fix commit was only added to row 29 after
# row 30 had already been written (matches original authoring order).
$ws.Range("F29").Value = "6369b468d2611907f95278f0fdac075f4a37b030"

# Give the new date cells (column A) the same date format / wrap style as the
# rest of the column by copying formatting down from the row above.
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A27:A30").PasteSpecial(-4122) | Out-Null

# Restore the row heights used elsewhere for rows whose remark wraps onto
# multiple lines.
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 60

$excel.CutCopyMode = $false

# Grow the worksheet table ("Tableau1") so the new rows participate in the
# filter / styling, mirroring the 4 additional journal entries.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F30")) | Out-Null

# Update the view so the freshly entered rows are visible/selected, matching
# where the author was working when they saved.
$ws.Range("E33").Select() | Out-Null

Write-Host "Journal entries for 2022-05-06 (serial 44687) added."
